$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

function Set-CellWithStyleFrom {
    param([string]$SourceAddress, [string]$TargetAddress, [double]$Value)
    $ws.Range($SourceAddress).Copy() | Out-Null
    $ws.Range($TargetAddress).PasteSpecial($xlPasteFormats) | Out-Null
    $ws.Range($TargetAddress).Value = $Value
}

# Row 9 (Болотная Виктория): add K9 = 5, style copied from J9
Set-CellWithStyleFrom "J9" "K9" 5

# Row 12 (Гарас Кристина): J12 gets value 5 (style unchanged), add K12 = 5 (style copied from I12)
$ws.Range("J12").Value = 5
Set-CellWithStyleFrom "I12" "K12" 5

# Row 16 (Дубкова Вероника): add J16 = 5 and K16 = 5 (style copied from J9, which is style "5")
Set-CellWithStyleFrom "J9" "J16" 5
Set-CellWithStyleFrom "J9" "K16" 5

# Row 20 (Ковшов Глеб): add J20 = 5 and K20 = 5 (style copied from I19, which is style "7")
Set-CellWithStyleFrom "I19" "J20" 5
Set-CellWithStyleFrom "I19" "K20" 5

# Row 22 (Кувшинова Ирина): add K22 = 5 (style copied from J22)
Set-CellWithStyleFrom "J22" "K22" 5

# Row 24 (Кунаева Кира): add K24 = 5 (style copied from J24)
Set-CellWithStyleFrom "J24" "K24" 5

# Row 26 (Миргасимов Расим): add K26 = 5 (style copied from J26)
Set-CellWithStyleFrom "J26" "K26" 5

# Row 28 (Решетняк Денис): add K28 = 5 (style copied from J28)
Set-CellWithStyleFrom "J28" "K28" 5

# Update the view: active selection in the bottom-right frozen pane becomes K20
$ws.Range("K20").Select() | Out-Null

Write-Host "Edit complete"
